$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

# H2 is never edited by this update and keeps the original 'General' style (s=3);
# use it as a format donor so forcing text on percent-looking values (which Excel
# would otherwise auto-convert to a numeric percentage) doesn't leave the cell on a
# different style index (quote-prefixed) than before.
$ws.Range("H2").Copy() | Out-Null

$ws.Range("E2").Value = '2026-02-09 05:18:37'
$ws.Range("E3").Value = '2026-02-09 05:18:40'
$ws.Range("H3").Value = '''94%'
$ws.Range("H3").PasteSpecial(-4122) | Out-Null
$ws.Range("M3").Value = '-5.3 °C 4:56 TU'
$ws.Range("O3").Value = '-6.2 °C'
$ws.Range("E4").Value = '2026-02-09 05:18:43'
$ws.Range("I4").Value = '0.1 mm'
$ws.Range("J4").Value = '1008.2 hPa'
$ws.Range("K4").Value = '-0.1 MJ/m2'
$ws.Range("N4").Value = '2.9 °C 4:59 TU'
$ws.Range("O4").Value = '4.5 °C'
$ws.Range("E5").Value = '2026-02-09 05:18:45'
$ws.Range("L5").Value = '12.6 km/h - 166º 4:58 TU'
$ws.Range("M5").Value = '-4.8 °C 4:59 TU'
$ws.Range("E6").Value = '2026-02-09 05:18:48'
$ws.Range("N6").Value = '4.7 °C 4:59 TU'
$ws.Range("O6").Value = '6.6 °C'
$ws.Range("E7").Value = '2026-02-09 05:18:51'
$ws.Range("H7").Value = '''71%'
$ws.Range("H7").PasteSpecial(-4122) | Out-Null
$ws.Range("E8").Value = '2026-02-09 05:18:53'
$ws.Range("H8").Value = '''79%'
$ws.Range("H8").PasteSpecial(-4122) | Out-Null
$ws.Range("E9").Value = '2026-02-09 05:18:56'
$ws.Range("H9").Value = '''87%'
$ws.Range("H9").PasteSpecial(-4122) | Out-Null
$ws.Range("N9").Value = '3.0 °C 4:53 TU'
$ws.Range("O9").Value = '7.0 °C'
$ws.Range("E10").Value = '2026-02-09 05:18:59'
$ws.Range("N10").Value = '3.0 °C 4:55 TU'
$ws.Range("O10").Value = '5.7 °C'
$ws.Range("E11").Value = '2026-02-09 05:19:01'
$ws.Range("H11").Value = '''98%'
$ws.Range("H11").PasteSpecial(-4122) | Out-Null
$ws.Range("E12").Value = '2026-02-09 05:19:03'
$ws.Range("H12").Value = '''91%'
$ws.Range("H12").PasteSpecial(-4122) | Out-Null
$ws.Range("N12").Value = '4.6 °C 4:58 TU'
$ws.Range("O12").Value = '7.5 °C'
$ws.Range("E13").Value = '2026-02-09 05:19:06'
$ws.Range("K13").Value = '-0.1 MJ/m2'
$ws.Range("O13").Value = '-1.7 °C'
$ws.Range("E14").Value = '2026-02-09 05:19:09'
$ws.Range("O14").Value = '7.4 °C'
$ws.Range("E15").Value = '2026-02-09 05:19:11'
$ws.Range("H15").Value = '''87%'
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$ws.Range("O15").Value = '5.6 °C'
$ws.Range("E16").Value = '2026-02-09 05:19:14'
$ws.Range("L16").Value = '20.9 km/h - 228º 4:58 TU'
$ws.Range("E17").Value = '2026-02-09 05:19:17'
$ws.Range("H17").Value = '''95%'
$ws.Range("H17").PasteSpecial(-4122) | Out-Null
$ws.Range("E18").Value = '2026-02-09 05:19:20'
$ws.Range("N18").Value = '3.8 °C 4:45 TU'
$ws.Range("O18").Value = '6.3 °C'
$ws.Range("E19").Value = '2026-02-09 05:19:23'
$ws.Range("N19").Value = '2.8 °C 4:58 TU'
$ws.Range("E20").Value = '2026-02-09 05:19:25'
$ws.Range("E21").Value = '2026-02-09 05:19:28'
$ws.Range("J21").Value = '1009.9 hPa'
$ws.Range("N21").Value = '-0.9 °C 4:50 TU'
$ws.Range("O21").Value = '0.6 °C'
$ws.Range("E22").Value = '2026-02-09 05:19:31'
$ws.Range("E23").Value = '2026-02-09 05:19:33'
$ws.Range("H23").Value = '''83%'
$ws.Range("H23").PasteSpecial(-4122) | Out-Null
$ws.Range("E24").Value = '2026-02-09 05:19:36'
$ws.Range("H24").Value = '''88%'
$ws.Range("H24").PasteSpecial(-4122) | Out-Null
$ws.Range("O24").Value = '4.3 °C'
$ws.Range("E25").Value = '2026-02-09 05:19:39'
$ws.Range("H25").Value = '''76%'
$ws.Range("H25").PasteSpecial(-4122) | Out-Null
$ws.Range("E26").Value = '2026-02-09 05:19:41'
$ws.Range("E27").Value = '2026-02-09 05:19:44'
$ws.Range("K27").Value = '-0.1 MJ/m2'
$ws.Range("N27").Value = '-4.6 °C 4:43 TU'
$ws.Range("O27").Value = '-4.2 °C'
$ws.Range("E28").Value = '2026-02-09 05:19:47'
$ws.Range("H28").Value = '''96%'
$ws.Range("H28").PasteSpecial(-4122) | Out-Null
$ws.Range("O28").Value = '3.6 °C'
$ws.Range("E29").Value = '2026-02-09 05:19:50'
$ws.Range("E30").Value = '2026-02-09 05:19:52'
$ws.Range("N30").Value = '4.8 °C 4:36 TU'
$ws.Range("O30").Value = '6.6 °C'
$ws.Range("E31").Value = '2026-02-09 05:19:55'
$ws.Range("K31").Value = '-0.1 MJ/m2'
$ws.Range("E32").Value = '2026-02-09 05:19:57'
$ws.Range("H32").Value = '''80%'
$ws.Range("H32").PasteSpecial(-4122) | Out-Null
$ws.Range("E33").Value = '2026-02-09 05:20:00'
$ws.Range("J33").Value = '1009.8 hPa'
$ws.Range("N33").Value = '-2.0 °C 4:49 TU'
$ws.Range("O33").Value = '-0.6 °C'
$ws.Range("E34").Value = '2026-02-09 05:20:03'
$ws.Range("O34").Value = '-3.2 °C'
$ws.Range("E35").Value = '2026-02-09 05:20:06'
$ws.Range("E36").Value = '2026-02-09 05:20:08'
$ws.Range("H36").Value = '''82%'
$ws.Range("H36").PasteSpecial(-4122) | Out-Null
$ws.Range("N36").Value = '5.4 °C 4:59 TU'
$ws.Range("O36").Value = '8.5 °C'
$ws.Range("E37").Value = '2026-02-09 05:20:11'
$ws.Range("L37").Value = '14.4 km/h - 20º 4:59 TU'
$ws.Range("N37").Value = '1.2 °C 4:59 TU'
$ws.Range("O37").Value = '3.4 °C'
$ws.Range("E38").Value = '2026-02-09 05:20:13'
$ws.Range("H38").Value = '''99%'
$ws.Range("H38").PasteSpecial(-4122) | Out-Null
$ws.Range("N38").Value = '4.4 °C 4:59 TU'
$ws.Range("O38").Value = '6.1 °C'
$ws.Range("E39").Value = '2026-02-09 05:20:16'
$ws.Range("H39").Value = '''83%'
$ws.Range("H39").PasteSpecial(-4122) | Out-Null
$ws.Range("E40").Value = '2026-02-09 05:20:18'
$ws.Range("E41").Value = '2026-02-09 05:20:21'
$ws.Range("H41").Value = '''56%'
$ws.Range("H41").PasteSpecial(-4122) | Out-Null
$ws.Range("J41").Value = '1008.1 hPa'
$ws.Range("L41").Value = '23.0 km/h - 289º 4:42 TU'
$ws.Range("O41").Value = '10.9 °C'
$ws.Range("E42").Value = '2026-02-09 05:20:23'
$ws.Range("H42").Value = '''97%'
$ws.Range("H42").PasteSpecial(-4122) | Out-Null
$ws.Range("O42").Value = '6.5 °C'
$ws.Range("E43").Value = '2026-02-09 05:20:26'
$ws.Range("N43").Value = '5.9 °C 4:59 TU'
$ws.Range("O43").Value = '6.4 °C'
$ws.Range("E44").Value = '2026-02-09 05:20:29'
$ws.Range("O44").Value = '-7.4 °C'
$ws.Range("E45").Value = '2026-02-09 05:20:31'
$ws.Range("E46").Value = '2026-02-09 05:20:34'
$ws.Range("O46").Value = '6.2 °C'
